# Gene_temp_data.xlsx — updates following Ann's comments to paper draft 1 11/28/20
#
# Rename the two rate-category labels used in the "Rate_Data" sheet's
# Rate_Type column (C2:C145):
#   "Constitutive_Rate" -> "Microbe_Independent_Rate"
#   "Induced_Rate"      -> "Microbe_Dependent_Rate"
# and leave the selection on C7 (the last interactive cell touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rate_Data")

$rateRange = $ws.Range("C2:C145")

$rateRange.Replace("Constitutive_Rate", "Microbe_Independent_Rate", 1) | Out-Null
$rateRange.Replace("Induced_Rate", "Microbe_Dependent_Rate", 1) | Out-Null

$ws.Activate() | Out-Null
$ws.Range("C7").Select() | Out-Null
